$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 1367.5  # H29
$ws.Cells.Item(29, 10).Value = 3500  # J29
$ws.Cells.Item(29, 12).Value = 10500  # L29
$ws.Cells.Item(29, 14).Value = -11062  # N29

$ws.Cells.Item(32, 8).Value = 9500  # H32
$ws.Cells.Item(32, 9).Value = 0  # I32
$ws.Cells.Item(32, 11).Value = 0  # K32
$ws.Cells.Item(32, 13).ClearContents()  # M32

$ws.Cells.Item(39, 8).Value = 562.1875  # H39
$ws.Cells.Item(39, 10).Value = 694.5  # J39
$ws.Cells.Item(39, 12).Value = 2083.5  # L39
$ws.Cells.Item(39, 14).Value = -2675.5  # N39

$ws.Cells.Item(53, 8).Value = 288.5  # H53
$ws.Cells.Item(53, 9).Value = 511.75  # I53
$ws.Cells.Item(53, 10).Value = 65.25  # J53
$ws.Cells.Item(53, 11).Value = 511.75  # K53
$ws.Cells.Item(53, 12).Value = 65.25  # L53
$ws.Cells.Item(53, 13).Value = 125.25  # M53
$ws.Cells.Item(53, 14).Value = -1339.25  # N53

$ws.Cells.Item(62, 8).Value = 5481.3237  # H62
$ws.Cells.Item(62, 9).Value = 4332.0386  # I62
$ws.Cells.Item(62, 10).Value = 9216.5  # J62
$ws.Cells.Item(62, 11).Value = 4332.0386  # K62
$ws.Cells.Item(62, 12).Value = 9216.5  # L62
$ws.Cells.Item(62, 13).Value = -3708.0386  # M62
$ws.Cells.Item(62, 14).Value = -10464.5  # N62

$ws.Cells.Item(64, 8).Value = 4820  # H64
$ws.Cells.Item(64, 9).Value = 4730  # I64
$ws.Cells.Item(64, 11).Value = 4730  # K64
$ws.Cells.Item(64, 13).Value = -4482  # M64

$ws.Cells.Item(65, 8).Value = 5481.3237  # H65
$ws.Cells.Item(65, 9).Value = 4332.0386  # I65
$ws.Cells.Item(65, 10).Value = 9216.5  # J65
$ws.Cells.Item(65, 11).Value = 21660.193  # K65
$ws.Cells.Item(65, 12).Value = 46082.5  # L65
$ws.Cells.Item(65, 13).Value = -18540.193  # M65
$ws.Cells.Item(65, 14).Value = -52322.5  # N65

$ws.Cells.Item(67, 8).Value = 4820  # H67
$ws.Cells.Item(67, 9).Value = 4730  # I67
$ws.Cells.Item(67, 11).Value = 4730  # K67
$ws.Cells.Item(67, 13).Value = -3872  # M67

$ws.Cells.Item(74, 8).Value = 12985.571  # H74
$ws.Cells.Item(74, 10).Value = 17450  # J74
$ws.Cells.Item(74, 12).Value = 17450  # L74
$ws.Cells.Item(74, 14).Value = -19322  # N74

$ws.Cells.Item(76, 8).Value = 8314  # H76
$ws.Cells.Item(76, 9).Value = 0  # I76
$ws.Cells.Item(76, 10).Value = 8314  # J76
$ws.Cells.Item(76, 11).Value = 0  # K76
$ws.Cells.Item(76, 12).Value = 8314  # L76
$ws.Cells.Item(76, 13).ClearContents()  # M76
$ws.Cells.Item(76, 14).Value = -8944  # N76

$ws.Cells.Item(77, 8).Value = 12985.571  # H77
$ws.Cells.Item(77, 10).Value = 17450  # J77
$ws.Cells.Item(77, 12).Value = 87250  # L77
$ws.Cells.Item(77, 14).Value = -96610  # N77

$ws.Cells.Item(79, 8).Value = 8314  # H79
$ws.Cells.Item(79, 9).Value = 0  # I79
$ws.Cells.Item(79, 10).Value = 8314  # J79
$ws.Cells.Item(79, 11).Value = 0  # K79
$ws.Cells.Item(79, 12).Value = 8314  # L79
$ws.Cells.Item(79, 13).ClearContents()  # M79
$ws.Cells.Item(79, 14).Value = -10498  # N79

$ws.Cells.Item(80, 8).Value = 908.4  # H80
$ws.Cells.Item(80, 10).Value = 848.5714  # J80
$ws.Cells.Item(80, 12).Value = 2545.7142  # L80
$ws.Cells.Item(80, 14).Value = -4541.7142  # N80

$ws.Cells.Item(83, 8).Value = 908.4  # H83
$ws.Cells.Item(83, 10).Value = 848.5714  # J83
$ws.Cells.Item(83, 12).Value = 7637.1426  # L83
$ws.Cells.Item(83, 14).Value = -17621.1426  # N83

$ws.Cells.Item(86, 8).Value = 3187.6843  # H86
$ws.Cells.Item(86, 9).Value = 1904.7142  # I86
$ws.Cells.Item(86, 10).Value = 6780  # J86
$ws.Cells.Item(86, 11).Value = 1904.7142  # K86
$ws.Cells.Item(86, 12).Value = 6780  # L86
$ws.Cells.Item(86, 13).Value = -781.7141999999999  # M86
$ws.Cells.Item(86, 14).Value = -9026  # N86

$ws.Cells.Item(89, 8).Value = 3187.6843  # H89
$ws.Cells.Item(89, 9).Value = 1904.7142  # I89
$ws.Cells.Item(89, 10).Value = 6780  # J89
$ws.Cells.Item(89, 11).Value = 9523.571  # K89
$ws.Cells.Item(89, 12).Value = 33900  # L89
$ws.Cells.Item(89, 13).Value = -3907.571  # M89
$ws.Cells.Item(89, 14).Value = -45132  # N89

$ws.Cells.Item(96, 8).Value = 222.73914  # H96
$ws.Cells.Item(96, 9).Value = 223.86667  # I96
$ws.Cells.Item(96, 10).Value = 220.625  # J96
$ws.Cells.Item(96, 11).Value = 671.60001  # K96
$ws.Cells.Item(96, 12).Value = 661.875  # L96
$ws.Cells.Item(96, 13).Value = 701.39999  # M96
$ws.Cells.Item(96, 14).Value = -3407.875  # N96

$ws.Cells.Item(98, 8).Value = 1668144.6  # H98
$ws.Cells.Item(98, 9).Value = 1895405.5  # I98
$ws.Cells.Item(98, 11).Value = 1895405.5  # K98
$ws.Cells.Item(98, 13).Value = -1893907.5  # M98

$ws.Cells.Item(111, 8).Value = 9595758  # H111
$ws.Cells.Item(111, 9).Value = 14287250  # I111
$ws.Cells.Item(111, 11).Value = 42861750  # K111
$ws.Cells.Item(111, 13).Value = -42858683  # M111

$ws.Cells.Item(112, 8).Value = 3327.361  # H112
$ws.Cells.Item(112, 9).Value = 1482.75  # I112
$ws.Cells.Item(112, 10).Value = 3854.3928  # J112
$ws.Cells.Item(112, 11).Value = 4448.25  # K112
$ws.Cells.Item(112, 12).Value = 11563.1784  # L112
$ws.Cells.Item(112, 13).Value = -3340.25  # M112
$ws.Cells.Item(112, 14).Value = -13779.1784  # N112

$ws.Cells.Item(116, 8).Value = 7846.8335  # H116
$ws.Cells.Item(116, 9).Value = 9035  # I116
$ws.Cells.Item(116, 10).Value = 6658.6665  # J116
$ws.Cells.Item(116, 11).Value = 9035  # K116
$ws.Cells.Item(116, 12).Value = 6658.6665  # L116
$ws.Cells.Item(116, 13).Value = -5593  # M116
$ws.Cells.Item(116, 14).Value = -13542.6665  # N116

$ws.Cells.Item(122, 8).Value = 1668144.6  # H122
$ws.Cells.Item(122, 9).Value = 1895405.5  # I122
$ws.Cells.Item(122, 11).Value = 5686216.5  # K122
$ws.Cells.Item(122, 13).Value = -5683766.5  # M122

$ws.Cells.Item(132, 8).Value = 2116.014  # H132
$ws.Cells.Item(132, 9).Value = 1981.8853  # I132
$ws.Cells.Item(132, 11).Value = 5945.6559  # K132
$ws.Cells.Item(132, 13).Value = -3415.6559  # M132

$ws.Cells.Item(137, 8).Value = 2564.5334  # H137
$ws.Cells.Item(137, 9).Value = 2173.238  # I137
$ws.Cells.Item(137, 10).Value = 3477.5557  # J137
$ws.Cells.Item(137, 11).Value = 6519.714  # K137
$ws.Cells.Item(137, 12).Value = 10432.6671  # L137
$ws.Cells.Item(137, 13).Value = -3969.714  # M137
$ws.Cells.Item(137, 14).Value = -15532.6671  # N137

$ws.Cells.Item(138, 8).Value = 2410.31  # H138
$ws.Cells.Item(138, 9).Value = 1094.125  # I138
$ws.Cells.Item(138, 10).Value = 3287.7666  # J138
$ws.Cells.Item(138, 11).Value = 3282.375  # K138
$ws.Cells.Item(138, 12).Value = 9863.299800000001  # L138
$ws.Cells.Item(138, 13).Value = 1857.625  # M138
$ws.Cells.Item(138, 14).Value = -20143.2998  # N138

$ws.Cells.Item(141, 8).Value = 2043.5883  # H141
$ws.Cells.Item(141, 9).Value = 1021.51166  # I141
$ws.Cells.Item(141, 11).Value = 3064.53498  # K141
$ws.Cells.Item(141, 13).Value = 2115.46502  # M141

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 240.22728  # H5
$ws.Cells.Item(5, 9).Value = 229.3  # I5
$ws.Cells.Item(5, 11).Value = 229.3  # K5
$ws.Cells.Item(5, 13).Value = -117.3  # M5

$ws.Cells.Item(31, 8).Value = 23867.75  # H31
$ws.Cells.Item(31, 9).Value = 13490.333  # I31
$ws.Cells.Item(31, 11).Value = 13490.333  # K31
$ws.Cells.Item(31, 13).Value = -13196.333  # M31

$ws.Cells.Item(32, 8).Value = 9808.152  # H32
$ws.Cells.Item(32, 9).Value = 4823.4473  # I32
$ws.Cells.Item(32, 11).Value = 4823.4473  # K32
$ws.Cells.Item(32, 13).Value = -4536.4473  # M32

$ws.Cells.Item(45, 8).Value = 84617850  # H45
$ws.Cells.Item(45, 9).Value = 122224296  # I45
$ws.Cells.Item(45, 10).Value = 3327.5  # J45
$ws.Cells.Item(45, 11).Value = 122224296  # K45
$ws.Cells.Item(45, 12).Value = 3327.5  # L45
$ws.Cells.Item(45, 13).Value = -122223919  # M45
$ws.Cells.Item(45, 14).Value = -4081.5  # N45

$ws.Cells.Item(61, 8).Value = 3045.1936  # H61
$ws.Cells.Item(61, 9).Value = 2499.5862  # I61
$ws.Cells.Item(61, 11).Value = 2499.5862  # K61
$ws.Cells.Item(61, 13).Value = -2287.5862  # M61

$ws.Cells.Item(63, 8).Value = 3702.1904  # H63
$ws.Cells.Item(63, 9).Value = 2731.5  # I63
$ws.Cells.Item(63, 11).Value = 2731.5  # K63
$ws.Cells.Item(63, 13).Value = -2045.5  # M63

$ws.Cells.Item(66, 8).Value = 3702.1904  # H66
$ws.Cells.Item(66, 9).Value = 2731.5  # I66
$ws.Cells.Item(66, 11).Value = 13657.5  # K66
$ws.Cells.Item(66, 13).Value = -10225.5  # M66

$ws.Cells.Item(74, 8).Value = 4575.7407  # H74
$ws.Cells.Item(74, 9).Value = 2537.7222  # I74
$ws.Cells.Item(74, 11).Value = 2537.7222  # K74
$ws.Cells.Item(74, 13).Value = -1663.7222  # M74

$ws.Cells.Item(77, 8).Value = 4575.7407  # H77
$ws.Cells.Item(77, 9).Value = 2537.7222  # I77
$ws.Cells.Item(77, 11).Value = 12688.611  # K77
$ws.Cells.Item(77, 13).Value = -8320.611000000001  # M77

$ws.Cells.Item(132, 8).Value = 4000.8333  # H132
$ws.Cells.Item(132, 9).Value = 2678.4194  # I132
$ws.Cells.Item(132, 10).Value = 12199.8  # J132
$ws.Cells.Item(132, 11).Value = 8035.2582  # K132
$ws.Cells.Item(132, 12).Value = 36599.39999999999  # L132
$ws.Cells.Item(132, 13).Value = -5505.2582  # M132
$ws.Cells.Item(132, 14).Value = -41659.39999999999  # N132

$ws.Cells.Item(136, 8).Value = 3045.1936  # H136
$ws.Cells.Item(136, 9).Value = 2499.5862  # I136
$ws.Cells.Item(136, 11).Value = 7498.758600000001  # K136
$ws.Cells.Item(136, 13).Value = -4948.758600000001  # M136

$ws.Cells.Item(139, 8).Value = 101884.445  # H139
$ws.Cells.Item(139, 9).Value = 0  # I139
$ws.Cells.Item(139, 10).Value = 101884.445  # J139
$ws.Cells.Item(139, 11).Value = 0  # K139
$ws.Cells.Item(139, 12).Value = 101884.445  # L139
$ws.Cells.Item(139, 13).ClearContents()  # M139
$ws.Cells.Item(139, 14).Value = -112164.445  # N139

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 240.22728  # H4
$ws.Cells.Item(4, 9).Value = 229.3  # I4
$ws.Cells.Item(4, 11).Value = 229.3  # K4
$ws.Cells.Item(4, 13).Value = -114.3  # M4

$ws.Cells.Item(37, 8).Value = 7993.143  # H37
$ws.Cells.Item(37, 9).Value = 1988  # I37
$ws.Cells.Item(37, 11).Value = 1988  # K37
$ws.Cells.Item(37, 13).Value = -1851  # M37

$ws.Cells.Item(56, 8).Value = 0  # H56
$ws.Cells.Item(56, 10).Value = 0  # J56
$ws.Cells.Item(56, 12).Value = 0  # L56
$ws.Cells.Item(56, 14).ClearContents()  # N56

$ws.Cells.Item(86, 8).Value = 7629.087  # H86
$ws.Cells.Item(86, 9).Value = 4571.636  # I86
$ws.Cells.Item(86, 10).Value = 10431.75  # J86
$ws.Cells.Item(86, 11).Value = 4571.636  # K86
$ws.Cells.Item(86, 12).Value = 10431.75  # L86
$ws.Cells.Item(86, 13).Value = -3448.636  # M86
$ws.Cells.Item(86, 14).Value = -12677.75  # N86

$ws.Cells.Item(89, 8).Value = 7629.087  # H89
$ws.Cells.Item(89, 9).Value = 4571.636  # I89
$ws.Cells.Item(89, 10).Value = 10431.75  # J89
$ws.Cells.Item(89, 11).Value = 22858.18  # K89
$ws.Cells.Item(89, 12).Value = 52158.75  # L89
$ws.Cells.Item(89, 13).Value = -17242.18  # M89
$ws.Cells.Item(89, 14).Value = -63390.75  # N89

$ws.Cells.Item(94, 8).Value = 768.46155  # H94
$ws.Cells.Item(94, 9).Value = 416  # I94
$ws.Cells.Item(94, 10).Value = 4998  # J94
$ws.Cells.Item(94, 11).Value = 416  # K94
$ws.Cells.Item(94, 12).Value = 4998  # L94
$ws.Cells.Item(94, 13).Value = 35  # M94
$ws.Cells.Item(94, 14).Value = -5900  # N94

$ws.Cells.Item(99, 8).Value = 1270  # H99
$ws.Cells.Item(99, 9).Value = 1270  # I99
$ws.Cells.Item(99, 11).Value = 1270  # K99
$ws.Cells.Item(99, 13).Value = 228  # M99

$ws.Cells.Item(107, 8).Value = 2383.2  # H107
$ws.Cells.Item(107, 9).Value = 2216.5715  # I107
$ws.Cells.Item(107, 10).Value = 2529  # J107
$ws.Cells.Item(107, 11).Value = 2216.5715  # K107
$ws.Cells.Item(107, 12).Value = 2529  # L107
$ws.Cells.Item(107, 13).Value = -296.5715  # M107
$ws.Cells.Item(107, 14).Value = -6369  # N107

$ws.Cells.Item(134, 8).Value = 4772.706  # H134
$ws.Cells.Item(134, 9).Value = 3703.7693  # I134
$ws.Cells.Item(134, 11).Value = 11111.3079  # K134
$ws.Cells.Item(134, 13).Value = -8576.3079  # M134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 2151.9092  # H22
$ws.Cells.Item(22, 9).Value = 3437.3333  # I22
$ws.Cells.Item(22, 10).Value = 609.4  # J22
$ws.Cells.Item(22, 11).Value = 3437.3333  # K22
$ws.Cells.Item(22, 12).Value = 609.4  # L22
$ws.Cells.Item(22, 13).Value = -3087.3333  # M22
$ws.Cells.Item(22, 14).Value = -1309.4  # N22

$ws.Cells.Item(31, 8).Value = 5586.1777  # H31
$ws.Cells.Item(31, 9).Value = 2192.5925  # I31
$ws.Cells.Item(31, 10).Value = 10676.556  # J31
$ws.Cells.Item(31, 11).Value = 2192.5925  # K31
$ws.Cells.Item(31, 12).Value = 10676.556  # L31
$ws.Cells.Item(31, 13).Value = -1897.5925  # M31
$ws.Cells.Item(31, 14).Value = -11266.556  # N31

$ws.Cells.Item(34, 8).Value = 5586.1777  # H34
$ws.Cells.Item(34, 9).Value = 2192.5925  # I34
$ws.Cells.Item(34, 10).Value = 10676.556  # J34
$ws.Cells.Item(34, 11).Value = 2192.5925  # K34
$ws.Cells.Item(34, 12).Value = 10676.556  # L34
$ws.Cells.Item(34, 13).Value = -1990.5925  # M34
$ws.Cells.Item(34, 14).Value = -11080.556  # N34

$ws.Cells.Item(111, 8).Value = 30000  # H111
$ws.Cells.Item(111, 10).Value = 30000  # J111
$ws.Cells.Item(111, 12).Value = 30000  # L111
$ws.Cells.Item(111, 14).Value = -38180  # N111

$ws.Cells.Item(132, 8).Value = 2672.4468  # H132
$ws.Cells.Item(132, 9).Value = 2151.4358  # I132
$ws.Cells.Item(132, 10).Value = 5212.375  # J132
$ws.Cells.Item(132, 11).Value = 6454.307400000001  # K132
$ws.Cells.Item(132, 12).Value = 15637.125  # L132
$ws.Cells.Item(132, 13).Value = -3924.307400000001  # M132
$ws.Cells.Item(132, 14).Value = -20697.125  # N132

$ws.Cells.Item(134, 8).Value = 3133.6316  # H134
$ws.Cells.Item(134, 9).Value = 2419.7144  # I134
$ws.Cells.Item(134, 10).Value = 5132.6  # J134
$ws.Cells.Item(134, 11).Value = 7259.1432  # K134
$ws.Cells.Item(134, 12).Value = 15397.8  # L134
$ws.Cells.Item(134, 13).Value = -4724.1432  # M134
$ws.Cells.Item(134, 14).Value = -20467.8  # N134

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 4282416  # H4
$ws.Cells.Item(4, 10).Value = 51001  # J4
$ws.Cells.Item(4, 12).Value = 153003  # L4
$ws.Cells.Item(4, 14).Value = -153227  # N4

$ws.Cells.Item(75, 8).Value = 17641.666  # H75
$ws.Cells.Item(75, 9).Value = 1000  # I75
$ws.Cells.Item(75, 10).Value = 20970  # J75
$ws.Cells.Item(75, 11).Value = 3000  # K75
$ws.Cells.Item(75, 12).Value = 62910  # L75
$ws.Cells.Item(75, 13).Value = -2002  # M75
$ws.Cells.Item(75, 14).Value = -64906  # N75

$ws.Cells.Item(76, 8).Value = 18465.916  # H76
$ws.Cells.Item(76, 9).Value = 17329.5  # I76
$ws.Cells.Item(76, 10).Value = 19602.334  # J76
$ws.Cells.Item(76, 11).Value = 51988.5  # K76
$ws.Cells.Item(76, 12).Value = 58807.00199999999  # L76
$ws.Cells.Item(76, 13).Value = -51605.5  # M76
$ws.Cells.Item(76, 14).Value = -59573.00199999999  # N76

$ws.Cells.Item(78, 8).Value = 17641.666  # H78
$ws.Cells.Item(78, 9).Value = 1000  # I78
$ws.Cells.Item(78, 10).Value = 20970  # J78
$ws.Cells.Item(78, 11).Value = 9000  # K78
$ws.Cells.Item(78, 12).Value = 188730  # L78
$ws.Cells.Item(78, 13).Value = -4008  # M78
$ws.Cells.Item(78, 14).Value = -198714  # N78

$ws.Cells.Item(79, 8).Value = 18465.916  # H79
$ws.Cells.Item(79, 9).Value = 17329.5  # I79
$ws.Cells.Item(79, 10).Value = 19602.334  # J79
$ws.Cells.Item(79, 11).Value = 51988.5  # K79
$ws.Cells.Item(79, 12).Value = 58807.00199999999  # L79
$ws.Cells.Item(79, 13).Value = -50662.5  # M79
$ws.Cells.Item(79, 14).Value = -61459.00199999999  # N79

$ws.Cells.Item(116, 8).Value = 2569.7  # H116
$ws.Cells.Item(116, 9).Value = 1528.1428  # I116
$ws.Cells.Item(116, 10).Value = 5000  # J116
$ws.Cells.Item(116, 11).Value = 4584.428400000001  # K116
$ws.Cells.Item(116, 12).Value = 15000  # L116
$ws.Cells.Item(116, 13).Value = -1142.428400000001  # M116
$ws.Cells.Item(116, 14).Value = -21884  # N116

$ws.Cells.Item(122, 8).Value = 1139.6  # H122
$ws.Cells.Item(122, 9).Value = 1046.5  # I122
$ws.Cells.Item(122, 10).Value = 1201.6666  # J122
$ws.Cells.Item(122, 11).Value = 9418.5  # K122
$ws.Cells.Item(122, 12).Value = 10814.9994  # L122
$ws.Cells.Item(122, 13).Value = -6968.5  # M122
$ws.Cells.Item(122, 14).Value = -15714.9994  # N122

$ws.Cells.Item(134, 8).Value = 2970.1667  # H134
$ws.Cells.Item(134, 9).Value = 2282.5557  # I134
$ws.Cells.Item(134, 10).Value = 5033  # J134
$ws.Cells.Item(134, 11).Value = 6847.6671  # K134
$ws.Cells.Item(134, 12).Value = 15099  # L134
$ws.Cells.Item(134, 13).Value = -1777.6671  # M134
$ws.Cells.Item(134, 14).Value = -25239  # N134

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2106  # H80
$ws.Cells.Item(80, 9).Value = 2275.6667  # I80
$ws.Cells.Item(80, 10).Value = 1902.4  # J80
$ws.Cells.Item(80, 11).Value = 2275.6667  # K80
$ws.Cells.Item(80, 12).Value = 1902.4  # L80
$ws.Cells.Item(80, 13).Value = -1277.6667  # M80
$ws.Cells.Item(80, 14).Value = -3898.4  # N80

$ws.Cells.Item(83, 8).Value = 2106  # H83
$ws.Cells.Item(83, 9).Value = 2275.6667  # I83
$ws.Cells.Item(83, 10).Value = 1902.4  # J83
$ws.Cells.Item(83, 11).Value = 11378.3335  # K83
$ws.Cells.Item(83, 12).Value = 9512  # L83
$ws.Cells.Item(83, 13).Value = -6386.333500000001  # M83
$ws.Cells.Item(83, 14).Value = -19496  # N83

$ws.Cells.Item(97, 8).Value = 815.3125  # H97
$ws.Cells.Item(97, 9).Value = 509.44446  # I97
$ws.Cells.Item(97, 11).Value = 509.44446  # K97
$ws.Cells.Item(97, 13).Value = -13.44445999999999  # M97

$ws.Cells.Item(102, 8).Value = 3254.8975  # H102
$ws.Cells.Item(102, 9).Value = 2148.7307  # I102
$ws.Cells.Item(102, 11).Value = 2148.7307  # K102
$ws.Cells.Item(102, 13).Value = -526.7307000000001  # M102

$ws.Cells.Item(107, 8).Value = 15152006  # H107
$ws.Cells.Item(107, 9).Value = 23809752  # I107
$ws.Cells.Item(107, 11).Value = 23809752  # K107
$ws.Cells.Item(107, 13).Value = -23807832  # M107

$ws.Cells.Item(123, 8).Value = 44718.4  # H123
$ws.Cells.Item(123, 10).Value = 44718.4  # J123
$ws.Cells.Item(123, 12).Value = 44718.4  # L123
$ws.Cells.Item(123, 14).Value = -49618.4  # N123

$ws.Cells.Item(126, 8).Value = 3404.6829  # H126
$ws.Cells.Item(126, 9).Value = 3292.0527  # I126
$ws.Cells.Item(126, 11).Value = 9876.158100000001  # K126
$ws.Cells.Item(126, 13).Value = -7406.158100000001  # M126

$ws.Cells.Item(132, 8).Value = 3826.8362  # H132
$ws.Cells.Item(132, 9).Value = 3439.1353  # I132
$ws.Cells.Item(132, 10).Value = 4424.5415  # J132
$ws.Cells.Item(132, 11).Value = 10317.4059  # K132
$ws.Cells.Item(132, 12).Value = 13273.6245  # L132
$ws.Cells.Item(132, 13).Value = -7787.4059  # M132
$ws.Cells.Item(132, 14).Value = -18333.6245  # N132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2521.05  # H46
$ws.Cells.Item(46, 9).Value = 684  # I46
$ws.Cells.Item(46, 10).Value = 2980.3125  # J46
$ws.Cells.Item(46, 11).Value = 684  # K46
$ws.Cells.Item(46, 12).Value = 2980.3125  # L46
$ws.Cells.Item(46, 13).Value = -496  # M46
$ws.Cells.Item(46, 14).Value = -3356.3125  # N46

$ws.Cells.Item(51, 8).Value = 0  # H51
$ws.Cells.Item(51, 10).Value = 0  # J51
$ws.Cells.Item(51, 12).Value = 0  # L51
$ws.Cells.Item(51, 14).ClearContents()  # N51

$ws.Cells.Item(55, 8).Value = 716.5862  # H55
$ws.Cells.Item(55, 9).Value = 191.85  # I55
$ws.Cells.Item(55, 10).Value = 1882.6666  # J55
$ws.Cells.Item(55, 11).Value = 191.85  # K55
$ws.Cells.Item(55, 12).Value = 1882.6666  # L55
$ws.Cells.Item(55, 13).Value = -18.84999999999999  # M55
$ws.Cells.Item(55, 14).Value = -2228.6666  # N55

$ws.Cells.Item(57, 8).Value = 26999.5  # H57
$ws.Cells.Item(57, 10).Value = 24000  # J57
$ws.Cells.Item(57, 12).Value = 24000  # L57
$ws.Cells.Item(57, 14).Value = -25132  # N57

$ws.Cells.Item(61, 8).Value = 1916.5143  # H61
$ws.Cells.Item(61, 9).Value = 1531.7059  # I61
$ws.Cells.Item(61, 10).Value = 15000  # J61
$ws.Cells.Item(61, 11).Value = 1531.7059  # K61
$ws.Cells.Item(61, 12).Value = 15000  # L61
$ws.Cells.Item(61, 13).Value = -1329.7059  # M61
$ws.Cells.Item(61, 14).Value = -15404  # N61

$ws.Cells.Item(99, 8).Value = 62071.25  # H99

$ws.Cells.Item(113, 8).Value = 1916.5143  # H113
$ws.Cells.Item(113, 9).Value = 1531.7059  # I113
$ws.Cells.Item(113, 10).Value = 15000  # J113
$ws.Cells.Item(113, 11).Value = 1531.7059  # K113
$ws.Cells.Item(113, 12).Value = 15000  # L113
$ws.Cells.Item(113, 13).Value = 638.2941000000001  # M113
$ws.Cells.Item(113, 14).Value = -19340  # N113

$ws.Cells.Item(122, 8).Value = 5264.45  # H122
$ws.Cells.Item(122, 9).Value = 5143.3  # I122
$ws.Cells.Item(122, 11).Value = 15429.9  # K122
$ws.Cells.Item(122, 13).Value = -12979.9  # M122

$ws.Cells.Item(132, 8).Value = 3020.2144  # H132
$ws.Cells.Item(132, 9).Value = 1588.1666  # I132
$ws.Cells.Item(132, 10).Value = 5597.9  # J132
$ws.Cells.Item(132, 11).Value = 4764.4998  # K132
$ws.Cells.Item(132, 12).Value = 16793.7  # L132
$ws.Cells.Item(132, 13).Value = -2234.4998  # M132
$ws.Cells.Item(132, 14).Value = -21853.7  # N132

$ws.Cells.Item(136, 8).Value = 8200.096  # H136
$ws.Cells.Item(136, 9).Value = 4613.684  # I136
$ws.Cells.Item(136, 10).Value = 9461.981  # J136
$ws.Cells.Item(136, 11).Value = 13841.052  # K136
$ws.Cells.Item(136, 12).Value = 28385.943  # L136
$ws.Cells.Item(136, 13).Value = -11291.052  # M136
$ws.Cells.Item(136, 14).Value = -33485.943  # N136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 793453.9  # H14
$ws.Cells.Item(14, 9).Value = 27263.637  # I14
$ws.Cells.Item(14, 10).Value = 5007500  # J14
$ws.Cells.Item(14, 11).Value = 27263.637  # K14
$ws.Cells.Item(14, 12).Value = 5007500  # L14
$ws.Cells.Item(14, 13).Value = -27095.637  # M14
$ws.Cells.Item(14, 14).Value = -5007836  # N14

$ws.Cells.Item(132, 8).Value = 2697.8447  # H132
$ws.Cells.Item(132, 9).Value = 2377.725  # I132
$ws.Cells.Item(132, 10).Value = 3409.2222  # J132
$ws.Cells.Item(132, 11).Value = 7133.174999999999  # K132
$ws.Cells.Item(132, 12).Value = 10227.6666  # L132
$ws.Cells.Item(132, 13).Value = -4603.174999999999  # M132
$ws.Cells.Item(132, 14).Value = -15287.6666  # N132
